# Updated Masterdata as per 2nd may Data Refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A "regulatory center" id corrections - rows where 10002 should be 10003
$ws.Range("A3").Value = 10003
$ws.Range("A23").Value = 10003
$ws.Range("A43").Value = 10003
$ws.Range("A63").Value = 10003
$ws.Range("A83").Value = 10003

# Column A "regulatory center" id corrections - rows where 10005 should be 10003
$ws.Range("A105").Value = 10003
$ws.Range("A114").Value = 10003
$ws.Range("A123").Value = 10003
$ws.Range("A132").Value = 10003
$ws.Range("A141").Value = 10003

# Refreshed view state: selection moved past the data block (row 162 onward)
$ws.Rows("162:1048576").Select()
